$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# STEP 1: Insert a new "Meta description" paragraph right after the
# document title (Heading1) paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Review of Aloha Fruit Bonanza slot: gameplay, payouts, free spins. Play for free and win up to 7,500x the stake. RTP of 97.01%.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# STEP 2: At the end of the document there are two trailing paragraphs:
#   - a duplicate bold title ("Play Aloha Fruit Bonanza Free: ...")
#   - an italic meta-description ("Review of Aloha Fruit Bonanza ...")
# Drop the duplicate bold title paragraph entirely and replace the
# italic paragraph's text with the new image-generation prompt, keeping
# the italic run formatting.
# ---------------------------------------------------------------------
$boldTitlePara = $null
$italicDescPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text.StartsWith("Play Aloha Fruit Bonanza Free")) {
        $boldTitlePara = $para
    }
    if ($text.StartsWith("Review of Aloha Fruit Bonanza slot")) {
        $italicDescPara = $para
    }
}

$tailRange = $d.Range($boldTitlePara.Range.Start, $italicDescPara.Range.End)

$promptXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Prompt: Create a cartoon-style feature image for the game "Aloha Fruit Bonanza". The image should feature a happy Maya warrior wearing glasses. The design should have a tropical feel, with bright colors and a beach background. It should also prominently feature fruit symbols from the game, such as watermelon, coconut cocktail, and dragon fruit cocktail, as well as the red number seven. The warrior should be holding a slot machine lever, and there should be cascading symbols falling around him. The overall feeling of the image should be fun and exciting, capturing the lightheartedness of the game.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $tailRange.InsertXML($promptXml)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
